$wb = $excel.ActiveWorkbook

# --- Sprint2 sheet: mark US17 (row 6) and US18 (row 5) as Done, fill in
#     actual size/time and completion date ---
$sprint2 = $wb.Worksheets.Item("Sprint2")

$sprint2.Range("D5").Value = "Done"
$sprint2.Range("G5").Value = 15
$sprint2.Range("H5").Value = 60
$sprint2.Range("I2").Copy()
$sprint2.Range("I5").PasteSpecial(-4122)  # xlPasteFormats
$sprint2.Range("I5").Value = 43186

$sprint2.Range("D6").Value = "Done"
$sprint2.Range("G6").Value = 30
$sprint2.Range("H6").Value = 60
$sprint2.Range("I2").Copy()
$sprint2.Range("I6").PasteSpecial(-4122)  # xlPasteFormats
$sprint2.Range("I6").Value = 43186

# --- Stories sheet: highlight the completed stories (US01-US11) ---
$stories = $wb.Worksheets.Item("Stories")
$stories.Range("A2:A11").Interior.Color = 65535

# --- Make Sprint2 the active tab/sheet ---
$sprint2.Activate()
[void]$sprint2.Range("L12").Select()
